$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.163103818893433
$ws.Range("B1").Value = 2.368344068527222
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.395913362503052
$ws.Range("E1").Value = 1.217627167701721
